# "Generate Report for Handoff"
#
# The localization status report is regenerated: the row that used to
# describe c1964d01-....md (which was "In Translation") has now been
# handed off, and the row that used to describe 1abd33cf-....md is now
# "Ready for handoff" with fresh handoff artifacts. The two rows swap
# display order (c1964d01 now first, 1abd33cf now second) on every
# sheet, and a couple of columns get a bit wider to fit the new text.

$wb = $excel.ActiveWorkbook

$url1abd = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5bdbc17bba90650aef396fa900cfb01ff0e3df1e/e2e/1abd33cf-6722-42c6-a51f-63d1fd08905b.md"
$urlc196 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5bdbc17bba90650aef396fa900cfb01ff0e3df1e/e2e/c1964d01-4667-4ad8-9491-9eb306304fcb.md"

# Closest achievable width to the target 17.2159881591797 model-width
# units given this engine's column-width quantization.
$newColWidth = 16.3

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "c1964d01-4667-4ad8-9491-9eb306304fcb.md"
$wsOverview.Range("B2").Value = "e2e\c1964d01-4667-4ad8-9491-9eb306304fcb.md"

$wsOverview.Range("A3").Value = "1abd33cf-6722-42c6-a51f-63d1fd08905b.md"
$wsOverview.Range("B3").Value = "e2e\1abd33cf-6722-42c6-a51f-63d1fd08905b.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-22 16:14:41"

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $urlc196, [System.Type]::Missing, [System.Type]::Missing, "e2e\c1964d01-4667-4ad8-9491-9eb306304fcb.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $url1abd, [System.Type]::Missing, [System.Type]::Missing, "e2e\1abd33cf-6722-42c6-a51f-63d1fd08905b.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "c1964d01-4667-4ad8-9491-9eb306304fcb.md"
$wsZhCn.Range("G2").Value = "c1964d01-4667-4ad8-9491-9eb306304fcb.eda0e898c828dd79cdce564f87a1185711c4f249.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "1abd33cf-6722-42c6-a51f-63d1fd08905b.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("G3").Value = "1abd33cf-6722-42c6-a51f-63d1fd08905b.458d0b96a15760efdc253cc90f55890ec1261b0e.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-22 16:14:36"

$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlc196, [System.Type]::Missing, [System.Type]::Missing, "c1964d01-4667-4ad8-9491-9eb306304fcb.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $url1abd, [System.Type]::Missing, [System.Type]::Missing, "1abd33cf-6722-42c6-a51f-63d1fd08905b.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "c1964d01-4667-4ad8-9491-9eb306304fcb.md"
$wsDeDe.Range("G2").Value = "c1964d01-4667-4ad8-9491-9eb306304fcb.eda0e898c828dd79cdce564f87a1185711c4f249.de-de.xlf"

$wsDeDe.Range("A3").Value = "1abd33cf-6722-42c6-a51f-63d1fd08905b.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("G3").Value = "1abd33cf-6722-42c6-a51f-63d1fd08905b.458d0b96a15760efdc253cc90f55890ec1261b0e.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-22 16:14:41"

$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlc196, [System.Type]::Missing, [System.Type]::Missing, "c1964d01-4667-4ad8-9491-9eb306304fcb.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $url1abd, [System.Type]::Missing, [System.Type]::Missing, "1abd33cf-6722-42c6-a51f-63d1fd08905b.md")
